$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AUTHOR_NAME (row 2) - fill in the author's name
$ws.Range("B2").Value = "Paul Ebbers"

# DN (row 4) - fill in the drawing number
$ws.Range("B4").Value = "210-05"

# FC-DATE (row 6) - today's date as a live formula, formatted as a short date.
# Set the number format first so Excel doesn't auto-apply its own implicit
# date format when the formula is entered into a "General" cell.
$ws.Range("B6").NumberFormat = "mm-dd-yy"
$ws.Range("B6").Formula = "=TODAY()"

# C6 picks up a date/time number format too (still center aligned, inherited
# from the existing style)
$ws.Range("C6").NumberFormat = "m/d/yy h:mm"

# Widen column B so the new values fit comfortably
$ws.Columns.Item(2).ColumnWidth = 25.6

# Leave the selection on C6, matching where the edit session ended up
$ws.Range("C6").Select()
